{"js": "// Replace the date line and every \"a\u00f7b=c, d\" answer cell with the new\n// values from the commit. Each old string is unique in the document, so a\n// plain body.search + insertText(\"Replace\") round-trip is safe and avoids\n// any positional/table-index assumptions.\nconst replacements = [\n  [\"2024-11-16 Saturday\", \"2024-11-17 Sunday\"],\n  [\"541\u00f78=67, 5\", \"781\u00f73=260, 1\"],\n  [\"784\u00f75=156, 4\", \"737\u00f77=105, 2\"],\n  [\"660\u00f79=73, 3\", \"663\u00f79=73, 6\"],\n  [\"664\u00f76=110, 4\", \"441\u00f79=49, 0\"],\n  [\"101\u00f77=14, 3\", \"469\u00f77=67, 0\"],\n  [\"781\u00f74=195, 1\", \"708\u00f77=101, 1\"],\n  [\"517\u00f75=103, 2\", \"862\u00f77=123, 1\"],\n  [\"352\u00f75=70, 2\", \"173\u00f76=28, 5\"],\n  [\"594\u00f77=84, 6\", \"271\u00f78=33, 7\"],\n  [\"587\u00f79=65, 2\", \"567\u00f77=81, 0\"],\n  [\"923\u00f74=230, 3\", \"360\u00f73=120, 0\"],\n  [\"211\u00f77=30, 1\", \"835\u00f79=92, 7\"],\n  [\"183\u00f74=45, 3\", \"833\u00f79=92, 5\"],\n  [\"923\u00f75=184, 3\", \"749\u00f78=93, 5\"],\n  [\"556\u00f76=92, 4\", \"105\u00f79=11, 6\"],\n  [\"208\u00f79=23, 1\", \"507\u00f76=84, 3\"],\n  [\"194\u00f79=21, 5\", \"384\u00f72=192, 0\"],\n  [\"740\u00f73=246, 2\", \"554\u00f74=138, 2\"],\n  [\"381\u00f77=54, 3\", \"860\u00f74=215, 0\"],\n  [\"903\u00f75=180, 3\", \"484\u00f74=121, 0\"],\n  [\"459\u00f72=229, 1\", \"794\u00f76=132, 2\"],\n  [\"269\u00f73=89, 2\", \"716\u00f76=119, 2\"],\n  [\"196\u00f72=98, 0\", \"712\u00f73=237, 1\"],\n  [\"159\u00f72=79, 1\", \"769\u00f73=256, 1\"],\n  [\"336\u00f74=84, 0\", \"179\u00f78=22, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00f7b=c, d\" answer cell with the new\n# values from the commit. Each old string is unique in the document, so a\n# Find/Replace over the whole story is safe and avoids any positional /\n# table-index assumptions.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-16 Saturday\", \"2024-11-17 Sunday\"),\n  @(\"541\u00f78=67, 5\", \"781\u00f73=260, 1\"),\n  @(\"784\u00f75=156, 4\", \"737\u00f77=105, 2\"),\n  @(\"660\u00f79=73, 3\", \"663\u00f79=73, 6\"),\n  @(\"664\u00f76=110, 4\", \"441\u00f79=49, 0\"),\n  @(\"101\u00f77=14, 3\", \"469\u00f77=67, 0\"),\n  @(\"781\u00f74=195, 1\", \"708\u00f77=101, 1\"),\n  @(\"517\u00f75=103, 2\", \"862\u00f77=123, 1\"),\n  @(\"352\u00f75=70, 2\", \"173\u00f76=28, 5\"),\n  @(\"594\u00f77=84, 6\", \"271\u00f78=33, 7\"),\n  @(\"587\u00f79=65, 2\", \"567\u00f77=81, 0\"),\n  @(\"923\u00f74=230, 3\", \"360\u00f73=120, 0\"),\n  @(\"211\u00f77=30, 1\", \"835\u00f79=92, 7\"),\n  @(\"183\u00f74=45, 3\", \"833\u00f79=92, 5\"),\n  @(\"923\u00f75=184, 3\", \"749\u00f78=93, 5\"),\n  @(\"556\u00f76=92, 4\", \"105\u00f79=11, 6\"),\n  @(\"208\u00f79=23, 1\", \"507\u00f76=84, 3\"),\n  @(\"194\u00f79=21, 5\", \"384\u00f72=192, 0\"),\n  @(\"740\u00f73=246, 2\", \"554\u00f74=138, 2\"),\n  @(\"381\u00f77=54, 3\", \"860\u00f74=215, 0\"),\n  @(\"903\u00f75=180, 3\", \"484\u00f74=121, 0\"),\n  @(\"459\u00f72=229, 1\", \"794\u00f76=132, 2\"),\n  @(\"269\u00f73=89, 2\", \"716\u00f76=119, 2\"),\n  @(\"196\u00f72=98, 0\", \"712\u00f73=237, 1\"),\n  @(\"159\u00f72=79, 1\", \"769\u00f73=256, 1\"),\n  @(\"336\u00f74=84, 0\", \"179\u00f78=22, 3\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $old\"\n  }\n}\n"}
